# Generate Report for Handback
# Updates the localization-status workbook: marks files as handed back and
# populates the "Latest Target File" / "Latest Handback File" columns (E/F)
# with hyperlinks, and refreshes the handback timestamps (column G).

$wb = $excel.ActiveWorkbook

$hyperlinkColor = 15570276  # OLE (BGR) value for RGB 0x6495ED (cornflower blue)

# ----------------------------------------------------------------------
# Overview sheet - shares the "Ready for handoff" / "Handed back..." string
# with the language sheets, so it must be refreshed too.
# ----------------------------------------------------------------------
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("B2").Value = "Handed back: in sync with en-US"
$wsOverview.Range("C2").Value = "Handed back: in sync with en-US"
$wsOverview.Range("B3").Value = "Handed back: in sync with en-US"
$wsOverview.Range("C3").Value = "Handed back: in sync with en-US"

function Set-HandbackHyperlink($ws, $cellRef, $displayText, $url) {
    $cell = $ws.Range($cellRef)
    $cell.Value = $displayText
    $ws.Hyperlinks.Add($cell, $url, "", "", $displayText)
    $cell.Font.Underline = 2
    $cell.Font.Color = $hyperlinkColor
}

# ----------------------------------------------------------------------
# zh-cn sheet
# ----------------------------------------------------------------------
$wsZh = $wb.Worksheets.Item("zh-cn")

$mdUrlZh = "https://github.com/OpenLocalizationTest/oltest/blob/03707ae99c51c269d015034c2baac24d8d3e877a/e2e/82cd0e9d-47d3-40fb-9dab-c63504f5cbb6.md"
$xlfUrlZh = "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/6b01366a7b29dadd2de61246f3878f071fb558a7/ol-handoff/OpenLocalizationTestOrg/oltest.zh-cn/yuwzho/82cd0e9d-47d3-40fb-9dab-c63504f5cbb6.1d1877cdb56ca16a4c696336e6b3429dbb4f778f.zh-cn.xlf"

# Status column now reflects the handback
$wsZh.Range("B2").Value = "Handed back: in sync with en-US"
$wsZh.Range("B3").Value = "Handed back: in sync with en-US"

# Latest Target File / Latest Handback File hyperlinks (columns E/F)
Set-HandbackHyperlink $wsZh "E2" "82cd0e9d-47d3-40fb-9dab-c63504f5cbb6.md" $mdUrlZh
Set-HandbackHyperlink $wsZh "F2" "82cd0e9d-47d3-40fb-9dab-c63504f5cbb6.1d1877cdb56ca16a4c696336e6b3429dbb4f778f.zh-cn.xlf" $xlfUrlZh
Set-HandbackHyperlink $wsZh "E3" "82cd0e9d-47d3-40fb-9dab-c63504f5cbb6.md" $mdUrlZh
Set-HandbackHyperlink $wsZh "F3" "82cd0e9d-47d3-40fb-9dab-c63504f5cbb6.1d1877cdb56ca16a4c696336e6b3429dbb4f778f.zh-cn.xlf" $xlfUrlZh

# Latest Handback DateTime (column G) / Handoff Reason (column H, unchanged text)
$wsZh.Range("G2").Value = "2016-02-06 04:23:01"
$wsZh.Range("G3").Value = "2016-02-06 04:23:01"
$wsZh.Range("H2").Value = "Include"
$wsZh.Range("H3").Value = "Include"

# ----------------------------------------------------------------------
# de-de sheet
# ----------------------------------------------------------------------
$wsDe = $wb.Worksheets.Item("de-de")

$mdUrlDe = "https://github.com/OpenLocalizationTest/oltest/blob/03707ae99c51c269d015034c2baac24d8d3e877a/e2e/82cd0e9d-47d3-40fb-9dab-c63504f5cbb6.md"
$xlfUrlDe = "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/8b0bc6a4f602abc9644fb8bf87c2eacb85f6c6e8/ol-handoff/OpenLocalizationTestOrg/oltest.de-de/yuwzho/82cd0e9d-47d3-40fb-9dab-c63504f5cbb6.1d1877cdb56ca16a4c696336e6b3429dbb4f778f.de-de.xlf"

$wsDe.Range("B2").Value = "Handed back: in sync with en-US"
$wsDe.Range("B3").Value = "Handed back: in sync with en-US"

Set-HandbackHyperlink $wsDe "E2" "82cd0e9d-47d3-40fb-9dab-c63504f5cbb6.md" $mdUrlDe
Set-HandbackHyperlink $wsDe "F2" "82cd0e9d-47d3-40fb-9dab-c63504f5cbb6.1d1877cdb56ca16a4c696336e6b3429dbb4f778f.de-de.xlf" $xlfUrlDe
Set-HandbackHyperlink $wsDe "E3" "82cd0e9d-47d3-40fb-9dab-c63504f5cbb6.md" $mdUrlDe
Set-HandbackHyperlink $wsDe "F3" "82cd0e9d-47d3-40fb-9dab-c63504f5cbb6.1d1877cdb56ca16a4c696336e6b3429dbb4f778f.de-de.xlf" $xlfUrlDe

$wsDe.Range("G2").Value = "2016-02-06 04:23:19"
$wsDe.Range("G3").Value = "2016-02-06 04:23:19"
$wsDe.Range("H2").Value = "Include"
$wsDe.Range("H3").Value = "Include"

Write-Host "Handback report generated."
